$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (OutSample -> TestSample)
$ws.Range("E1").Value = "RMSE_TestSample"
$ws.Range("F1").Value = "R2_TestSample"
$ws.Range("G1").Value = "Adjusted_R2_TestSample"

# Update numeric values for rows 2-9, columns B, E, F, G, H

# Row 2 (Num_Factors = 5)
$ws.Range("B2").Value = 0.02463865175384928
$ws.Range("E2").Value = 0.03543189471917484
$ws.Range("F2").Value = 0.9928022233450564
$ws.Range("G2").Value = 0.9921357625436728
$ws.Range("H2").Value = 291.2687803182

# Row 3 (Num_Factors = 6)
$ws.Range("B3").Value = 0.01908051223776923
$ws.Range("E3").Value = 0.0275803700828647
$ws.Range("F3").Value = 0.9953448168563664
$ws.Range("G3").Value = 0.9948178149910495
$ws.Range("H3").Value = 346.7694955966112

# Row 4 (Num_Factors = 7)
$ws.Range("B4").Value = 0.01604666133594783
$ws.Range("E4").Value = 0.02086707775059352
$ws.Range("F4").Value = 0.9974402804563148
$ws.Range("G4").Value = 0.9970957028254341
$ws.Range("H4").Value = 392.133218325583

# Row 5 (Num_Factors = 8)
$ws.Range("B5").Value = 0.0111091205506193
$ws.Range("E5").Value = 0.01548468615517097
$ws.Range("F5").Value = 0.9984218393605436
$ws.Range("G5").Value = 0.9981742847504328
$ws.Range("H5").Value = 473.0631617584364

# Row 6 (Num_Factors = 9)
$ws.Range("B6").Value = 0.007475639668471302
$ws.Range("E6").Value = 0.01009298211407882
$ws.Range("F6").Value = 0.9995341107955057
$ws.Range("G6").Value = 0.9994502507386968
$ws.Range("H6").Value = 593.2170164517049

# Row 7 (Num_Factors = 10)
$ws.Range("B7").Value = 0.005785531382205611
$ws.Range("E7").Value = 0.00822693784767531
$ws.Range("F7").Value = 0.9997546658653153
$ws.Range("G7").Value = 0.9997045976745632
$ws.Range("H7").Value = 718.1476233602531

# Row 8 (Num_Factors = 11)
$ws.Range("B8").Value = 0.005785531382205611
$ws.Range("E8").Value = 0.00822693784767531
$ws.Range("F8").Value = 0.9997546658653153
$ws.Range("G8").Value = 0.99969844345945
$ws.Range("H8").Value = 718.1476233602531

# Row 9 (Num_Factors = 12)
$ws.Range("B9").Value = 0.005785531382205611
$ws.Range("E9").Value = 0.00822693784767531
$ws.Range("F9").Value = 0.9997546658653153
$ws.Range("G9").Value = 0.9996920273628426
$ws.Range("H9").Value = 718.1476233602531
